# Updated cryptos list on Fri Oct 18 14:52:53 UTC 2024 with GitHub Actions
# Refresh Price (col D) and Volume(1h) (col E) for rows 2-51 on Sheet1.
# NumberFormat is forced to text ("@") before writing any D-column value
# that looks like a plain decimal number (e.g. "1.00", "599.10") so Excel
# keeps it as literal text instead of silently coercing it to a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.468.95'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '2.642.78'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.10'
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.42'
$ws.Range("E6").Value = '  +1.88%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").Value = '2.641.00'
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("E10").Value = '  +11.96%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("E15").Value = '  +5.75%  '
$ws.Range("D16").Value = '3.118.63'
$ws.Range("E16").Value = '  +1.43%  '
$ws.Range("D17").Value = '68.350.16'
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").Value = '2.635.97'
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("E19").Value = '  +3.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '369.51'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.41'
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.25'
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +1.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.12'
$ws.Range("E25").Value = '  +8.57%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("D28").Value = '2.765.40'
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("E29").Value = '  +3.78%  '
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '572.72'
$ws.Range("E32").Value = '  +3.64%  '
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("E34").Value = '  +3.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  +4.97%  '
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.06'
$ws.Range("E38").Value = '  +2.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.16'
$ws.Range("E39").Value = '  +1.72%  '
$ws.Range("E40").Value = '  +5.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.368'
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.34'
$ws.Range("E42").Value = '  +2.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.65'
$ws.Range("E43").Value = '  +3.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.61'
$ws.Range("E44").Value = '  +5.01%  '
$ws.Range("D45").Value = '0.0₆0324'
$ws.Range("E45").Value = '  +9.42%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.47'
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '155.56'
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.71'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.92'
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.70'
$ws.Range("E51").Value = '  +0.51%  '
